# Append the latest Adafruit IO reading as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 52

# Timestamp / Feed Key are plain text - Value assignment keeps them as text.
$ws.Cells.Item($newRow, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($newRow, 2).Value = "temperature"

# Value column holds a numeric-looking reading ("25") that must stay text,
# matching the rest of the column - format as Text before typing it in,
# just like entering '25 by hand in Excel.
$ws.Cells.Item($newRow, 3).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).Value = "25"

# Latitude / Longitude / Elevation were not reported for this reading.
$ws.Cells.Item($newRow, 4).Value = "N/A"
$ws.Cells.Item($newRow, 5).Value = "N/A"
$ws.Cells.Item($newRow, 6).Value = "N/A"
